# edit.ps1 - apply "feat: add 2022-Q3 data" change
#
# Summary of the change:
#  1. A new worksheet "2022-Q3" is inserted right after "总计" and before
#     "2022-Q1" (i.e. becomes the 2nd sheet), holding a fund-holdings table
#     for the new quarter.
#  2. The "总计" (summary) sheet gets a new row for "2022-Q3" inserted at
#     row 2; the rest of the quarters cascade down one row, and a new last
#     row ("2020-Q4") is appended at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned before "2022-Q1"
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q1")
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Range($cols[$i] + "1")
    $cell.Value2 = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$data = @(
    @("016935", "景顺长城中证500指数增强C", "15.57", "93.89", "2.09", "0.3254", 4),
    @("000978", "景顺长城量化精选股票", "7.14", "93.64", "1.96", "0.1399", 5),
    @("008851", "景顺长城量化对冲策略三个月定期开放灵活配置混合", "2.96", "64.77", "1.36", "0.0403", 5),
    @("006682", "景顺长城中证500指数增强A", "0.00", "93.89", "2.09", 0, 4)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $vals = $data[$r]

    $acell = $q3.Range("A" + $row)
    $acell.Value2 = $r
    $acell.Font.Bold = $true
    $acell.Borders.LineStyle = 1
    $acell.HorizontalAlignment = -4108
    $acell.VerticalAlignment = -4160

    # B: fund code (text)
    $q3.Range("B" + $row).NumberFormat = "@"
    $q3.Range("B" + $row).Value2 = $vals[0]

    # C: fund name (text)
    $q3.Range("C" + $row).NumberFormat = "@"
    $q3.Range("C" + $row).Value2 = $vals[1]

    # D: fund size (text)
    $q3.Range("D" + $row).NumberFormat = "@"
    $q3.Range("D" + $row).Value2 = $vals[2]

    # E: total stock position (text)
    $q3.Range("E" + $row).NumberFormat = "@"
    $q3.Range("E" + $row).Value2 = $vals[3]

    # F: position ratio (text)
    $q3.Range("F" + $row).NumberFormat = "@"
    $q3.Range("F" + $row).Value2 = $vals[4]

    # G: held market value (亿元) - text, except the last row which is a
    # genuine 0 number in the source data
    if ($row -eq 5) {
        $q3.Range("G" + $row).Value2 = 0
    } else {
        $q3.Range("G" + $row).NumberFormat = "@"
        $q3.Range("G" + $row).Value2 = $vals[5]
    }

    # H: position rank (number)
    $q3.Range("H" + $row).Value2 = $vals[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: cascade existing rows down and
#    insert the new "2022-Q3" figures at the top of the data block.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$b2 = $summary.Range("B2").Value2
$c2 = $summary.Range("C2").Value2
$d2 = $summary.Range("D2").Value2
$b3 = $summary.Range("B3").Value2
$c3 = $summary.Range("C3").Value2
$d3 = $summary.Range("D3").Value2
$b4 = $summary.Range("B4").Value2
$c4 = $summary.Range("C4").Value2
$d4 = $summary.Range("D4").Value2
$b5 = $summary.Range("B5").Value2
$c5 = $summary.Range("C5").Value2
$d5 = $summary.Range("D5").Value2

$summary.Range("B2").Value2 = "2022-Q3"
$summary.Range("C2").Value2 = 4
$summary.Range("D2").Value2 = 0.51

$summary.Range("B3").Value2 = $b2
$summary.Range("C3").Value2 = $c2
$summary.Range("D3").Value2 = $d2

$summary.Range("B4").Value2 = $b3
$summary.Range("C4").Value2 = $c3
$summary.Range("D4").Value2 = $d3

$summary.Range("B5").Value2 = $b4
$summary.Range("C5").Value2 = $c4
$summary.Range("D5").Value2 = $d4

# New row 6 ("2020-Q4"), cloning A5's style for the A6 index cell
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Range("A6").Value2 = 4
$summary.Range("B6").Value2 = $b5
$summary.Range("C6").Value2 = $c5
$summary.Range("D6").Value2 = $d5

Write-Output "2022-Q3 sheet added and 总计 updated"
